# Day-ahead dispatch data update: refresh line parameters R (col C) and X (col D)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C9").Value = 0.001
$ws.Range("D2:D9").Value = 0.0001

# Reflect the last selected cell in the sheet view
$ws.Range("G14").Select()

$wb.Save()
